# FINFLUX-2544  Automating Nabkisan Sanity Scenario
#
# Adds a new "Modify Transaction1" worksheet (after "Transactions") that
# records the steps of the "modify transaction / submit / navigate to loan"
# automation scenario, mirroring the existing NewLoanInput-style
# name/value layout.

$wb = $excel.ActiveWorkbook

# Source sheet we borrow existing cell formatting from, so the new sheet's
# styles line up with the rest of the workbook instead of minting unrelated
# ones.
$src = $wb.Worksheets.Item("NewLoanInput")

# New sheet goes after the last existing sheet (i.e. after "Transactions").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Modify Transaction1"

# --- Values -----------------------------------------------------------
$newSheet.Range("A1").Value = "OverDueTillDate"
$newSheet.Range("B1").Value = 42038

$newSheet.Range("A2").Value = "clickonsubmit"
$newSheet.Range("B2").Value = "Submit"

$newSheet.Range("A3").Value = "NavigateToLoan"
$newSheet.Range("B3").Value = "navigate"

# --- Formatting ---------------------------------------------------------
# Column A: shaded "label" look used throughout NewLoanInput.
$rA = $newSheet.Range("A1:A3")
$src.Range("A1").Copy()
$rA.PasteSpecial(-4122)   # xlPasteFormats
$rA.Font.Name = "Calibri"

# B1: date value -> reuse the existing shaded/date style.
$src.Range("B15").Copy()
$newSheet.Range("B1").PasteSpecial(-4122)   # xlPasteFormats

# B2:B3: text values -> same shaded style family, without word-wrap.
$rB = $newSheet.Range("B2:B3")
$src.Range("B15").Copy()
$rB.PasteSpecial(-4122)   # xlPasteFormats
$rB.WrapText = $false

$newSheet.Columns.Item(1).ColumnWidth = 15.86

# Make the new sheet the active / selected one, matching the saved view.
[void]$newSheet.Range("D5").Select()
